$wb = $excel.ActiveWorkbook

# --- "Affichage" sheet: course code PHY2701 -> PHY2710 ---
$ws1 = $wb.Worksheets.Item("Affichage")
$ws1.Range("A13").Value = "PHY2710"

# --- "Candidatures" sheet: prefix bare course codes with "PHY" ---
$ws2 = $wb.Worksheets.Item("Candidatures")
$ws2.Range("D2").Value = "PHY1441-30, PHY1620-30, PHY1620-30, PHY1620-90, PHY1620-90, PHY1651-30"
$ws2.Range("H2").Value = "PHY1441-90, PHY1620-90"
$ws2.Range("D3").Value = "PHY1441-30, PHY2701, PHY2701, PHY2701, PHY2701, PHY2701"
$ws2.Range("H3").Value = "PHY1441-90, PHY2710"
$ws2.Range("D4").Value = "PHY2710, PHY2710, PHY2400, PHY2400"
$ws2.Range("H4").Value = "PHY2710, PHY2400"

# --- Window / selection state ---
$excel.ActiveWindow.TabRatio = 0.984

[void]$ws1.Range("B23").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
[void]$ws2.Range("K39").Select()
